# Commit: "new commit 21/3 1"
# - bump the absPath / revisionPtr are workbook-session metadata the COM
#   surface doesn't expose; what we CAN reproduce is the actual content
#   change: a new "Sheet3" appended after Sheet2, holding a single
#   contact record (Amar Kadam) with two mailto hyperlinks on the email
#   cell (A3/A4, duplicated) and one on the password-looking cell (A5),
#   mirroring the shape of the existing Sheet2 contact-card sheet.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the end of the tab strip -------------------
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$ws.Name = "Sheet3"

# --- plain text rows (quote-prefixed, like the sibling sheets) -------
$ws.Range("A1").Value = "'Amar"
$ws.Range("A2").Value = "'Kadam"

# --- email, duplicated on two rows, each its own mailto hyperlink ----
$ws.Range("A3").Value = "amar.kadam@openxcell.com"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:amar.kadam@openxcell.com")
$ws.Range("A3").Value = "'amar.kadam@openxcell.com"

$ws.Range("A4").Value = "amar.kadam@openxcell.com"
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:amar.kadam@openxcell.com")
$ws.Range("A4").Value = "'amar.kadam@openxcell.com"

# --- password-looking cell, also hyperlinked the same way ------------
$ws.Range("A5").Value = "Amar@1234"
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:Amar@1234")
$ws.Range("A5").Value = "'Amar@1234"

# --- column A sized to fit the longest entry (best-effort match of
#     the authored 23.66-ish character width) -------------------------
$ws.Columns.Item(1).ColumnWidth = 22.8

# --- leave the selection/active cell on the last-entered cell --------
$ws.Range("A5").Select() | Out-Null
